$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings are kept as text (matching source formatting)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

$ws.Range("D2").Value = "69.336.97"
$ws.Range("E2").Value = "  -2.20%  "
$ws.Range("D3").Value = "3.683.12"
$ws.Range("E3").Value = "  -2.85%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "683.71"
$ws.Range("E5").Value = "  -2.97%  "
$ws.Range("D6").Value = "160.15"
$ws.Range("E6").Value = "  -6.09%  "
$ws.Range("D7").Value = "3.681.36"
$ws.Range("E7").Value = "  -2.89%  "
$ws.Range("E8").Value = "  -0.26%  "
$ws.Range("E9").Value = "  -5.58%  "
$ws.Range("E10").Value = "  -8.67%  "
$ws.Range("D11").Value = "7.19"
$ws.Range("E11").Value = "  -2.83%  "
$ws.Range("D12").Value = "0.436"
$ws.Range("E12").Value = "  -8.73%  "
$ws.Range("E13").Value = "  -6.42%  "
$ws.Range("D14").Value = "4.301.82"
$ws.Range("E14").Value = "  -2.86%  "
$ws.Range("E15").Value = "  -9.63%  "
$ws.Range("D16").Value = "3.682.99"
$ws.Range("E16").Value = "  -2.98%  "
$ws.Range("D17").Value = "69.369.11"
$ws.Range("E17").Value = "  -2.28%  "
$ws.Range("E18").Value = "  -1.14%  "
$ws.Range("D19").Value = "15.83"
$ws.Range("E19").Value = "  -9.11%  "
$ws.Range("D20").Value = "6.42"
$ws.Range("E20").Value = "  -9.96%  "
$ws.Range("D21").Value = "471.44"
$ws.Range("E21").Value = "  -7.97%  "
$ws.Range("D22").Value = "9.91"
$ws.Range("E22").Value = "  -4.21%  "
$ws.Range("E23").Value = "  -8.67%  "
$ws.Range("D24").Value = "79.45"
$ws.Range("E24").Value = "  -4.81%  "
$ws.Range("D25").Value = "3.827.38"
$ws.Range("E25").Value = "  -3.21%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  -9.01%  "
$ws.Range("D28").Value = "11.03"
$ws.Range("E28").Value = "  -11.92%  "
$ws.Range("D29").Value = "9.08"
$ws.Range("E29").Value = "  -11.41%  "
$ws.Range("E30").Value = "  -10.16%  "
$ws.Range("E31").Value = "  -13.19%  "
$ws.Range("D32").Value = "6.68"
$ws.Range("E32").Value = "  -8.59%  "
$ws.Range("E33").Value = "  -9.27%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").Value = "26.71"
$ws.Range("E35").Value = "  -8.04%  "
$ws.Range("D36").Value = "3.650.01"
$ws.Range("E36").Value = "  -2.81%  "
$ws.Range("E37").Value = "  -7.60%  "
$ws.Range("E38").Value = "  -10.82%  "
$ws.Range("D39").Value = "6.14"
$ws.Range("E39").Value = "  -3.67%  "
$ws.Range("E40").Value = "  -5.23%  "
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("E42").Value = "  -9.95%  "
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("D44").Value = "0.945"
$ws.Range("E44").Value = "  -6.23%  "
$ws.Range("D45").Value = "165.02"
$ws.Range("E45").Value = "  -2.96%  "
$ws.Range("D46").Value = "47.89"
$ws.Range("E46").Value = "  -3.80%  "
$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").Value = "2.71"
$ws.Range("E47").Value = "  -15.74%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "28.76"
$ws.Range("E48").Value = "  -1.82%  "
$ws.Range("D49").Value = "1.30"
$ws.Range("E49").Value = "  -4.64%  "
$ws.Range("E50").Value = "  -4.09%  "
$ws.Range("E51").Value = "  -10.67%  "
